# Adds two new columns, I0 (column I) and IF (column J), to the sheet,
# mirroring the structure of the existing header/data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - reuse the same style as the other header cells (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Data values for rows 2-34: row number, I value, J value
$data = @(
    @(2, 7, 8),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 6, 6),
    @(6, 10, 10),
    @(7, 9, 9),
    @(8, 7, 7),
    @(9, 6, 7),
    @(10, 6, 7),
    @(11, 9, 9),
    @(12, 6, 7),
    @(13, 6, 6),
    @(14, 6, 6),
    @(15, 7, 7),
    @(16, 7, 7),
    @(17, 9, 9),
    @(18, 6, 7),
    @(19, 9, 9),
    @(20, 8, 8),
    @(21, 6, 6),
    @(22, 6, 7),
    @(23, 6, 7),
    @(24, 9, 9),
    @(25, 8, 8),
    @(26, 6, 6),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 5, 5),
    @(31, 7, 7),
    @(32, 7, 8),
    @(33, 9, 9),
    @(34, 7, 7)
)

foreach ($item in $data) {
    $r = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
